# Updated filepaths for stat testing:
# swap the A1/B1 header labels ("unet" <-> "expert") and move the
# selection to B2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a1 = $ws.Range("A1").Value()
$b1 = $ws.Range("B1").Value()

$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

$ws.Range("B2").Select()
